# "Generate Report for Handoff"
#
# Adds two new tracked files (two .png screenshots) alongside the existing
# .md file, and refreshes the handoff timestamps / generated filenames for
# the existing .md entry across all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldGuid = "78302f29-1153-4969-a374-1280d8adc01a"
$newGuid = "31673876-2087-46c0-a480-fb2112b1ea72"

$mdName       = "$newGuid.md"
$pngName1     = "9f2b9db4-58b9-475c-9e65-7d5157416a0f.png"
$pngName2     = "bc9acc81-2120-4ce7-8dbc-d05e15a2bddb.png"

$zhHash       = "42c3367dcb17677218a83721fb05abee8cbfecae"
$mdTarget     = "$newGuid.$zhHash.zh-cn.xlf"
$deTarget     = "$newGuid.$zhHash.de-de.xlf"
$png1Target   = "476f890b81eb4a8a3da450a62cc599dc7d40e1e9.png"
$png2Target   = "39fc4f6c7ae8e5a57687f7538cd3d1c1ff848205.png"

$overviewDate = "2016-49-19 14:49:45"
$zhDate       = "2016-03-19 14:49:42"
$deDate       = "2016-03-19 14:49:45"
$epoch        = "0001-01-01 00:00:00"
$readyStatus  = "Ready for handoff"
$includeReason    = "Include"
$dependencyReason = "IsDependency"
$dependencyFrom   = "e2e\$mdName"

$mdCommit   = "a1b2c3d4e5f6071829384756a1b2c3d4e5f60718"
$zhCommit   = "66c91c953059f5c2d88695a95c5a69ec6efbc0c4"
$deCommit   = "1f1fcddd994b1b47e4fcd5c2ab74f2194001bfb2"

function MdUrl($name) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$name"
}
function XlfUrl($locale, $name) {
    if ($locale -eq "zh-cn") {
        return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$name"
    }
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$name"
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Refresh the existing .md row (row 2): new generated guid + new handoff date.
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Range("A2").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), (MdUrl $mdName), "", "", $mdName) | Out-Null
$wsOverview.Range("B2").Value = $readyStatus
$wsOverview.Range("C2").Value = $readyStatus
$wsOverview.Range("D2").Value = $overviewDate

# New row 3: first screenshot.
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
$wsOverview.Range("D3").Value = $overviewDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), (MdUrl $pngName1), "", "", $pngName1) | Out-Null

# New row 4: second screenshot.
$wsOverview.Range("B4").Value = $readyStatus
$wsOverview.Range("C4").Value = $readyStatus
$wsOverview.Range("D4").Value = $overviewDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), (MdUrl $pngName2), "", "", $pngName2) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" / "de-de": both share the same column layout:
# Source File Name | File Extension | Status | Latest Handoff File |
# Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------
function FillLocaleSheet($ws, $locale, $handoffDate, $xlfTarget) {
    # Row 2: refresh the existing .md entry.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("B2").Hyperlinks.Delete()
    $ws.Range("D2").Hyperlinks.Delete()

    $ws.Range("A2").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("A2"), (MdUrl $mdName), "", "", $mdName) | Out-Null

    $ws.Range("B2").Value = ".md"
    $ws.Hyperlinks.Add($ws.Range("B2"), (MdUrl $mdName), "", "", ".md") | Out-Null

    $ws.Range("C2").Value = $readyStatus
    $ws.Range("D2").Value = $xlfTarget
    $ws.Hyperlinks.Add($ws.Range("D2"), (XlfUrl $locale $xlfTarget), "", "", $xlfTarget) | Out-Null

    $ws.Range("E2").Value = $handoffDate
    $ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H2").Value = $epoch
    $ws.Range("I2").Value = $includeReason

    # Row 3: first screenshot (depends on the .md file).
    $ws.Range("B3").Value = ".png"
    $ws.Hyperlinks.Add($ws.Range("A3"), (MdUrl $pngName1), "", "", $pngName1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), (MdUrl $pngName1), "", "", ".png") | Out-Null

    $ws.Range("C3").Value = $readyStatus
    $ws.Range("D3").Value = $png1Target
    $ws.Hyperlinks.Add($ws.Range("D3"), (XlfUrl $locale $png1Target), "", "", $png1Target) | Out-Null

    $ws.Range("E3").Value = $handoffDate
    $ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H3").Value = $epoch
    $ws.Range("I3").Value = $dependencyReason
    $ws.Range("J3").Value = $dependencyFrom

    # Row 4: second screenshot (depends on the .md file).
    $ws.Range("B4").Value = ".png"
    $ws.Hyperlinks.Add($ws.Range("A4"), (MdUrl $pngName2), "", "", $pngName2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B4"), (MdUrl $pngName2), "", "", ".png") | Out-Null

    $ws.Range("C4").Value = $readyStatus
    $ws.Range("D4").Value = $png2Target
    $ws.Hyperlinks.Add($ws.Range("D4"), (XlfUrl $locale $png2Target), "", "", $png2Target) | Out-Null

    $ws.Range("E4").Value = $handoffDate
    $ws.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H4").Value = $epoch
    $ws.Range("I4").Value = $dependencyReason
    $ws.Range("J4").Value = $dependencyFrom
}

$wsZh = $wb.Worksheets.Item("zh-cn")
FillLocaleSheet $wsZh "zh-cn" $zhDate $mdTarget

$wsDe = $wb.Worksheets.Item("de-de")
FillLocaleSheet $wsDe "de-de" $deDate $deTarget
